$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 87, pushing the existing rows 87-94 down to 88-95
$ws.Rows.Item(87).Insert()

# Populate the new row 87 with the latest weekly price record
$ws.Cells.Item(87, 1).Value = 4
$ws.Cells.Item(87, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(87, 3).Value = "Los Lagos"
$ws.Cells.Item(87, 4).Value = 45265
$ws.Cells.Item(87, 5).Value = 10
$ws.Cells.Item(87, 6).Value = 300000000
$ws.Cells.Item(87, 7).Value = "Espárragos"
$ws.Cells.Item(87, 8).Value = "Sin especificar"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 500
$ws.Cells.Item(87, 11).Value = 2000
$ws.Cells.Item(87, 12).Value = 2000
$ws.Cells.Item(87, 13).Value = 2000
$ws.Cells.Item(87, 14).Value = "$/kilo"
$ws.Cells.Item(87, 15).Value = "Provincia de Linares"
$ws.Cells.Item(87, 16).Value = 2000
$ws.Cells.Item(87, 17).Value = 1
$ws.Cells.Item(87, 18).Value = "Hortaliza"
